$wb = $excel.ActiveWorkbook

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = 0.04468251151616323
$wsLinear.Range("B3").Value = 0.03102332466603601
$wsLinear.Range("B4").Value = 1295.041088241587
$wsLinear.Range("B5").Value = "[1.0, 0.1250755947006343, -0.006220958035794838, -0.02179810679630137, -0.031443926729890805, -0.037930929562663335, 0.1226935481085313, 0.2885840932616539, 0.10301226696570974, -0.039531350905854405, -0.04954548228505662, -0.047381127539237435, -0.04893331894029557, 0.11177809186989698, 0.3005717703972358, 0.1357518968706738, -0.0418467850464528, -0.04364448859338815, -0.02847051159053469, -0.048529177153476884]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = 0.7960526171721864
$wsNonLinear.Range("B5").Value = 0.133604303999067
$wsNonLinear.Range("B6").Value = 1245.263572625828
$wsNonLinear.Range("B7").Value = 3.515922037452919
$wsNonLinear.Range("B8").Value = -0.218894189237738
$wsNonLinear.Range("B9").Value = 1334.305495654982
$wsNonLinear.Range("B10").Value = "[1.0, 0.1250680618342565, 0.0018539440310870795, -0.013001136826587526, -0.023230530999312656, -0.029455740570979754, 0.12225909724910684, 0.27648184230987344, 0.1024708061896736, -0.03173219154539911, -0.04207220093836153, -0.04093916574905838, -0.04087523861511522, 0.11113025422000913, 0.2909240146007749, 0.13405462017163144, -0.03214510635049886, -0.03810990976421786, -0.021784958296120243, -0.040602951735164206]"
